$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 ("I0") and J1 ("IF") --------------------------
# Copy the existing header formatting (bold font + border + center/top
# alignment, same as B1:H1) onto the two new header cells before writing
# their text, so the added columns look consistent with the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data for I2:J76 ---------------------------------------------------
$iValues = @(8,9,9,7,8,8,12,8,7,9,6,7,7,9,7,8,7,8,9,8,9,9,9,8,8,9,9,9,9,9,9,9,9,9,9,9,9,9,9,7,8,7,7,7,6,6,7,7,8,6,8,6,7,6,7,5,8,9,9,10,7,7,7,7,8,8,3,7,8,8,8,7,6,7,6)
$jValues = @(8,10,9,8,8,8,12,8,7,9,6,7,7,9,7,8,7,8,9,8,9,9,9,9,8,9,9,9,9,10,9,9,9,9,9,9,9,9,9,8,8,7,7,8,6,7,7,7,8,7,8,8,7,7,8,6,8,9,9,10,7,7,7,8,8,8,4,7,8,8,8,7,6,7,6)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
